# Insert a new weekly price-report row at row 62 (pushing the existing
# rows 62-139 down to 63-140) and populate it with the new observation.
#
# The new row shares every descriptive column (Mercado, Región, Codreg,
# Categoría, Variedad, Calidad, Unidad de comercialización, Origen,
# Kg o Unidades, Clasificación, ...) with the row that is about to occupy
# row 63 immediately after the insert (the former row 62), so the
# quickest reliable way to seed it is to insert the blank row and then
# copy that row's values down into the freshly inserted row before
# overwriting the handful of cells (Fecha, Volumen, Precio mínimo,
# Precio máximo, Precio promedio ponderado, Precio $/Kg) that actually
# differ for the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 62:139 down to 63:140, leaving a blank row 62 behind.
$ws.Rows.Item(62).Insert()

# Seed the new row 62 with the same descriptive data as the row right
# below it (the former row 62, now sitting at row 63).
$ws.Range("A63:R63").Copy($ws.Range("A62:R62"))

# Overwrite the cells that differ for this new weekly observation.
$ws.Range("D62").Value = 45159
$ws.Range("J62").Value = 250
$ws.Range("K62").Value = 1300
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = 1380
$ws.Range("P62").Value = 690
